$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Molex Minifit")

# Add new row 29: 20-position Micro-Fit 3.0mm TH right-angle header
$row = 29
$ws.Cells.Item($row, 1).Value  = 430452001
$ws.Cells.Item($row, 2).Value  = 20
$ws.Cells.Item($row, 3).Value  = "3.0mm"
$ws.Cells.Item($row, 4).Value  = "3.0mm"
$ws.Cells.Item($row, 5).Value  = 90
$ws.Cells.Item($row, 6).Value  = "No"
$ws.Cells.Item($row, 7).Value  = "Gold"
$ws.Cells.Item($row, 8).Value  = 430452001
$ws.Cells.Item($row, 9).Value  = "Molex Micro-Fit.SchLib"
$ws.Cells.Item($row, 10).Value = 430452001
$ws.Cells.Item($row, 11).Value = "Molex Micro-Fit.PcbLib"
$ws.Cells.Item($row, 12).Value = "WM7492-ND"

# Digikey link + hyperlink for the new part, matching the style already
# used by the other Digikey-link cells in column M.
$ws.Hyperlinks.Add($ws.Cells.Item($row, 13), "https://www.digikey.com.au/en/products/detail/molex/0430452001/3044584")
$ws.Cells.Item($row, 13).Style = $ws.Cells.Item(27, 13).Style

# Restore cursor/selection to match the saved view state.
$ws.Range("P12").Select()
